$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 data: "usb error" / "MESA,E03\n" (literal backslash-n, matching existing data)
$ws.Range("A9").Value = "usb error"
$ws.Range("B9").Value = "MESA,E03\n"

# Update selection to D9 (matches diff's sheetView selection)
$ws.Range("D9").Select()
